# foodlist.xlsx update:
#   - Sheet2 ("diet categories"): add "vegetables"/"fruit"/"other" summary columns (C:E)
#   - Sheet1 ("diet table"): add a new "foliage" column (AJ) and refresh the
#     yogurt/pellets/crispbread/foliage totals for each species row

$wb = $excel.ActiveWorkbook

# --- Sheet2: add vegetables / fruit / other columns (C, D, E) -----------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("C1").Value = "vegetables"
$ws2.Range("D1").Value = "fruit"
$ws2.Range("E1").Value = "other"

$ws2.Range("B2").Value = 3
$ws2.Range("C2").Value = 14
$ws2.Range("D2").Value = 6
$ws2.Range("E2").Value = 5

$ws2.Range("B3").Value = 4
$ws2.Range("C3").Value = 15
$ws2.Range("D3").Value = 6
$ws2.Range("E3").Value = 5

$ws2.Range("B4").Value = 4
$ws2.Range("C4").Value = 14
$ws2.Range("D4").Value = 6
$ws2.Range("E4").Value = 5

$ws2.Range("B5").Value = 7
$ws2.Range("C5").Value = 14
$ws2.Range("D5").Value = 1
$ws2.Range("E5").Value = 5

$ws2.Range("B6").Value = 4
$ws2.Range("C6").Value = 15
$ws2.Range("D6").Value = 6
$ws2.Range("E6").Value = 5

$ws2.Range("C3").Select() | Out-Null

# --- Sheet1: add foliage column (AJ) and update the per-species totals --
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("AJ1").Value = "foliage"

$ws1.Range("AG2").Value = 21
$ws1.Range("AH2").Value = 21
$ws1.Range("AI2").Value = 21
$ws1.Range("AJ2").Value = 21

$ws1.Range("AG3").Value = 6
$ws1.Range("AH3").Value = 6
$ws1.Range("AI3").Value = 6
$ws1.Range("AJ3").Value = 6

$ws1.Range("AG4").Value = 13
$ws1.Range("AH4").Value = 13
$ws1.Range("AI4").Value = 13
$ws1.Range("AJ4").Value = 13

$ws1.Range("AG5").Value = 6
$ws1.Range("AH5").Value = 6
$ws1.Range("AI5").Value = 6
$ws1.Range("AJ5").Value = 17

$ws1.Range("AG6").Value = 7
$ws1.Range("AH6").Value = 7
$ws1.Range("AI6").Value = 7
$ws1.Range("AJ6").Value = 14

$ws1.Activate()
$ws1.Range("AJ7").Select() | Out-Null
# best-effort: scroll the view so column H is leftmost (matches the
# original author's scrolled viewport); not all hosts persist this.
$excel.ActiveWindow.ScrollColumn = 8
